$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1 with the same style as the other headers (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new data columns I and J for rows 2-6
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 8

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8
